# Snippet for line-type shapes (#267)
# Adds two new rows to the "Snippets" table on the "Snippets" worksheet:
#   ShapeCollection | addLine | excel-shape-lines | addStraightLine
#   Shape           | line    | excel-shape-lines | arrowLine

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Snippets")
$lo = $ws.ListObjects.Item("Snippets")

$row1 = $lo.ListRows.Add()
$row1.Range.Cells.Item(1, 1).Value = "ShapeCollection"
$row1.Range.Cells.Item(1, 2).Value = "addLine"
$row1.Range.Cells.Item(1, 3).Value = "excel-shape-lines"
$row1.Range.Cells.Item(1, 4).Value = "addStraightLine"

$row2 = $lo.ListRows.Add()
$row2.Range.Cells.Item(1, 1).Value = "Shape"
$row2.Range.Cells.Item(1, 2).Value = "line"
$row2.Range.Cells.Item(1, 3).Value = "excel-shape-lines"
$row2.Range.Cells.Item(1, 4).Value = "arrowLine"

# Move the selection down to the first empty row below the growing table,
# matching where Excel would leave the cursor after entering this data.
$ws.Range("A170").Select()
